$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# ---------------------------------------------------------------------------
# New log entries for WEDNESDAY (11/23/16, serial 42697) and THURSDAY
# (11/24/16, serial 42698), appended after the existing data (which ends at
# row 818). Each new data row is cloned (values + formatting) from an
# existing row elsewhere in the sheet that already carries the exact same
# staff/task/building/room/instructions combination, then the date (and, for
# a couple of rows, the task-type / note text) is patched to the new value.
# This guarantees styles/number-formats match exactly and avoids any
# shared-string transcription errors.
# ---------------------------------------------------------------------------

function Copy-LogRow($srcRow, $dstRow) {
    $src = $ws.Range("A$srcRow`:F$srcRow")
    $dst = $ws.Range("A$dstRow`:F$dstRow")
    $src.Copy($dst)
}

$wedSerial = 42697
$thuSerial = 42698

# --- Wednesday header/section -------------------------------------------------
Copy-LogRow 5 822
$ws.Range("B822").Value2 = $ws.Range("B56").Value2   # "WEDNESDAY"

Copy-LogRow 149 823
$ws.Range("B823").Value2 = $wedSerial
$ws.Rows.Item(823).RowHeight = 60

Copy-LogRow 150 824
$ws.Range("B824").Value2 = $wedSerial
$ws.Rows.Item(824).RowHeight = 75

Copy-LogRow 151 825
$ws.Range("B825").Value2 = $wedSerial
$ws.Rows.Item(825).RowHeight = 75

Copy-LogRow 667 826
$ws.Range("B826").Value2 = $wedSerial

Copy-LogRow 153 827
$ws.Range("B827").Value2 = $wedSerial
$ws.Rows.Item(827).RowHeight = 120

Copy-LogRow 435 828
$ws.Range("B828").Value2 = $wedSerial
$ws.Rows.Item(828).RowHeight = 45

Copy-LogRow 263 829
$ws.Range("B829").Value2 = $wedSerial

Copy-LogRow 766 830
$ws.Range("B830").Value2 = $wedSerial

Copy-LogRow 339 831
$ws.Range("B831").Value2 = $wedSerial

Copy-LogRow 185 832
$ws.Range("B832").Value2 = $wedSerial
$ws.Range("C832").Value2 = "1820"
$ws.Rows.Item(832).RowHeight = 90

Copy-LogRow 186 833
$ws.Range("B833").Value2 = $wedSerial
$ws.Range("C833").Value2 = "1820"
$ws.Rows.Item(833).RowHeight = 30

Copy-LogRow 248 834
$ws.Range("B834").Value2 = $wedSerial
$ws.Range("C834").Value2 = "1820"
$ws.Rows.Item(834).RowHeight = 30

Copy-LogRow 665 835
$ws.Range("B835").Value2 = $wedSerial
$ws.Range("C835").Value2 = "1820"
$ws.Range("F835").Value2 = "PLEAS BE ON TIME - GO EARLY - GUEST PROF ENDING EARLY TODAY ANYTIME FROM 18:20 - 18:30 pm. "
$ws.Rows.Item(835).RowHeight = 30

# --- Thursday header/section ---------------------------------------------------
Copy-LogRow 5 839
$ws.Range("B839").Value2 = $ws.Range("B5").Value2   # "THURSDAY"

Copy-LogRow 429 840
$ws.Range("B840").Value2 = $thuSerial
$ws.Rows.Item(840).RowHeight = 45

Copy-LogRow 299 841
$ws.Range("B841").Value2 = $thuSerial

Copy-LogRow 263 842
$ws.Range("B842").Value2 = $thuSerial

Copy-LogRow 302 843
$ws.Range("B843").Value2 = $thuSerial

Copy-LogRow 300 844
$ws.Range("B844").Value2 = $thuSerial

# --- View state -----------------------------------------------------------
$ws.Range("F848").Select()
